$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed line-flow results (380 kV case) for rows 2-25, columns B,D,E,F,G,L,M,N
$updates = @{
    "B2" = 1.597911568692552
    "D2" = 0.1696348709340896
    "E2" = 0.9374858084730135
    "F2" = 3.325952304526936
    "G2" = 0.002461777939971148
    "L2" = 0.5472208321561425
    "M2" = 0.460722331563268
    "N2" = 1.888593783721817
    "B3" = 1.535953024462344
    "D3" = 0.160534539045031
    "E3" = 0.81929992625399
    "F3" = 3.125374034908049
    "G3" = 0.002473307158535375
    "L3" = 0.5006151004010917
    "M3" = 0.4342876932094626
    "N3" = 1.91742430161932
    "B4" = 1.498905002768197
    "D4" = 0.1551231442069252
    "E4" = 0.7467512689009652
    "F4" = 3.00533716137133
    "G4" = 0.002480731608325262
    "L4" = 0.4723511345251552
    "M4" = 0.418359298421322
    "N4" = 1.935930530285514
    "B5" = 1.484055363383305
    "D5" = 0.1529609945352206
    "E5" = 0.7171831177621755
    "F5" = 2.957182940658242
    "G5" = 0.002483844472743102
    "L5" = 0.4609191155875862
    "M5" = 0.4119431309296928
    "N5" = 1.943674111416533
    "B6" = 1.481604485374618
    "D6" = 0.1526045336949551
    "E6" = 0.7122728639407683
    "F6" = 2.949232330325771
    "G6" = 0.002484366650554877
    "L6" = 0.4590259386459934
    "M6" = 0.410882210523873
    "N6" = 1.944972134181192
    "B7" = 1.498703735414352
    "D7" = 0.1550938121140604
    "E7" = 0.7463525296841595
    "F7" = 3.004684681661104
    "G7" = 0.002480773235266058
    "L7" = 0.4721966146229022
    "M7" = 0.4182724666446234
    "N7" = 1.93603414454065
    "B8" = 1.576340298053481
    "D8" = 0.1664597143287807
    "E8" = 0.8967263772593981
    "F8" = 3.256131795924119
    "G8" = 0.002465681799291309
    "L8" = 0.5310763835323655
    "M8" = 0.4515440632064909
    "N8" = 1.898367596154044
    "B9" = 1.736600457416046
    "D9" = 0.1902068269251629
    "E9" = 1.192157046046731
    "F9" = 3.775031038438385
    "G9" = 0.002438806942802193
    "L9" = 0.6494677606027039
    "M9" = 0.5192545452796651
    "N9" = 1.830886757408706
    "B10" = 1.859418108100897
    "D10" = 0.2086300398603385
    "E10" = 1.410139161574534
    "F10" = 4.173549648084247
    "G10" = 0.002420689133679079
    "L10" = 0.7384346858843003
    "M10" = 0.5706009114448705
    "N10" = 1.785201933239453
    "B11" = 1.916436398097972
    "D11" = 0.2172428534877611
    "E11" = 1.509642034631241
    "F11" = 4.35894916377498
    "G11" = 0.002412793488152647
    "L11" = 0.7793845440262714
    "M11" = 0.5943291791393364
    "N11" = 1.765264775559968
    "B12" = 1.938196229177152
    "D12" = 0.2205392983707952
    "E12" = 1.547381482198404
    "F12" = 4.429775354331582
    "G12" = 0.002409852865876249
    "L12" = 0.7949636715964914
    "M12" = 0.6033695160008818
    "N12" = 1.757836762293406
    "B13" = 1.933502325713334
    "D13" = 0.2198277691603323
    "E13" = 1.539250779091731
    "F13" = 4.41449365896591
    "G13" = 0.002410483997410686
    "L13" = 0.791605154543447
    "M13" = 0.6014200461729473
    "N13" = 1.759431098874993
    "B14" = 1.91822320205074
    "D14" = 0.2175133440146908
    "E14" = 1.512745618919496
    "F14" = 4.364763498125058
    "G14" = 0.002412550576659521
    "L14" = 0.7806647785166092
    "M14" = 0.595071821871926
    "N14" = 1.764651226896085
    "B15" = 1.908886304059649
    "D15" = 0.2161002928306743
    "E15" = 1.496518565015663
    "F15" = 4.334383858323974
    "G15" = 0.002413822819588256
    "L15" = 0.7739730051602862
    "M15" = 0.5911905589096165
    "N15" = 1.767864569364569
    "B16" = 1.855715196153199
    "D16" = 0.2080719760998306
    "E16" = 1.403644199305432
    "F16" = 4.161518475137314
    "G16" = 0.002421212057053467
    "L16" = 0.7357684399169102
    "M16" = 0.5690578112571103
    "N16" = 1.786521898329688
    "B17" = 1.823392671222962
    "D17" = 0.2032074055690032
    "E17" = 1.346763816891837
    "F17" = 4.056543962637306
    "G17" = 0.002425833440267437
    "L17" = 0.7124561581429134
    "M17" = 0.5555761639363936
    "N17" = 1.798184174689366
    "B18" = 1.804909460469787
    "D18" = 0.2004311717130349
    "E18" = 1.314079223831584
    "F18" = 3.996550641298711
    "G18" = 0.002428524157395193
    "L18" = 0.6990923574671513
    "M18" = 0.547856650920366
    "N18" = 1.804971537352422
    "B19" = 1.798669791791383
    "D19" = 0.1994948723895504
    "E19" = 1.303017871562474
    "F19" = 3.976303255052358
    "G19" = 0.002429440805789794
    "L19" = 0.6945751947583005
    "M19" = 0.5452488789745757
    "N19" = 1.807283274250747
    "B20" = 1.826822274477649
    "D20" = 0.203722986151206
    "E20" = 1.352815494253292
    "F20" = 4.067678606331697
    "G20" = 0.002425338113866432
    "L20" = 0.7149331278851605
    "M20" = 0.5570076974231029
    "N20" = 1.796934472811735
    "B21" = 1.922706459071549
    "D21" = 0.2181921850072399
    "E21" = 1.520529114648241
    "F21" = 4.379353411553154
    "G21" = 0.002411942238595099
    "L21" = 0.7838762411394953
    "M21" = 0.5969349440732259
    "N21" = 1.763114642409707
    "B22" = 1.986354525974605
    "D22" = 0.2278531916422537
    "E22" = 1.630493891426539
    "F22" = 4.586675127360479
    "G22" = 0.002403474306274335
    "L22" = 0.8293577260946279
    "M22" = 0.6233508605637468
    "N22" = 1.741721232774054
    "B23" = 1.952293400555504
    "D23" = 0.2226776727271158
    "E23" = 1.571767558251025
    "F23" = 4.475682442917048
    "G23" = 0.002407967703948126
    "L23" = 0.8050434848280759
    "M23" = 0.6092222340152631
    "N23" = 1.753074270050238
    "B24" = 1.825271441353152
    "D24" = 0.2034898285226063
    "E24" = 1.350079479719398
    "F24" = 4.062643521185635
    "G24" = 0.002425561945386946
    "L24" = 0.7138131694606784
    "M24" = 0.5563604038740095
    "N24" = 1.797499205877656
    "B25" = 1.69236870692049
    "D25" = 0.183617615068016
    "E25" = 1.112116246639232
    "F25" = 3.631734593082939
    "G25" = 0.002445789323335357
    "L25" = 0.6171060244535624
    "M25" = 0.5006633674261138
    "N25" = 1.848458603566895
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
